$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New item rows 11-17 -----------------------------------------------

# Row 11: HK MP5A2 Plastic Fixed Stock (typed in place, keeps existing formatting)
$ws.Range("A11").Value = "hk_mp5a2_plastic_fixed_stock"
$ws.Range("B11").Value = "HK MP5A2 Plastic Fixed Stock"
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = 0.26
$ws.Range("E11").Value = -10
$ws.Range("F11").Value = -12
$ws.Range("M11").Value = 0

# Row 12: HK MP5A3 Early Generation Stock Endplate (pasted values, no formatting,
# unused cells dropped entirely)
$ws.Range("A12").Value = "hk_mp5a3_early_gen_stock_endplate"
$ws.Range("B12").Value = "HK MP5A3 Early Generation Stock Endplate"
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 0.08
$ws.Range("M12").Value = 1500
$ws.Range("A12:D12").ClearFormats()
$ws.Range("M12").ClearFormats()
$ws.Range("E12:L12").Clear()

# Row 13: HK MP5A3 Early Generation Collapsible Stock (pasted values)
$ws.Range("A13").Value = "hk_mp5a3_early_gen_collapsible_stock"
$ws.Range("B13").Value = "HK MP5A3 Early Generation Collapsible Stock"
$ws.Range("C13").Value = 9
$ws.Range("D13").Value = 0.12
$ws.Range("E13").Value = -9
$ws.Range("F13").Value = -8
$ws.Range("M13").Value = 0
$ws.Range("A13:F13").ClearFormats()
$ws.Range("M13").ClearFormats()
$ws.Range("G13:L13").Clear()

# Row 14: HP MP5/HK94 Stock Endcap (pasted values)
$ws.Range("A14").Value = "hk_mp5_hk94_stock_endcap"
$ws.Range("B14").Value = "HP MP5/HK94 Stock Endcap"
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 0.06
$ws.Range("M14").Value = 1000
$ws.Range("A14:D14").ClearFormats()
$ws.Range("M14").ClearFormats()
$ws.Range("E14:L14").Clear()

# Row 15: HK MP5/HK94 Choate Stock Base (typed in place, keeps formatting)
$ws.Range("A15").Value = "hk_mp5_hk94_choate_stock_base"
$ws.Range("B15").Value = "HK MP5/HK94 Choate Stock Base"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0.06
$ws.Range("M15").Value = 750

# Row 16: HK MP5/HK94 Choate Stock (typed in place, keeps formatting)
$ws.Range("A16").Value = "hk_mp5_hk94_choate_stock"
$ws.Range("B16").Value = "HK MP5/HK94 Choate Stock"
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 0.12
$ws.Range("E16").Value = -8
$ws.Range("F16").Value = -13
$ws.Range("M16").Value = 0

# Row 17: HK MP5 Endcap Sling Swivel (typed in place, keeps formatting)
$ws.Range("A17").Value = "hk_mp5_endcap_sling_swivel"
$ws.Range("B17").Value = "HK MP5 Endcap Sling Swivel"
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 0.01
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 10
$ws.Range("H17").Value = 0.2
$ws.Range("M17").Value = 250

# --- Extend the shared "balance" formula from N4:N9 down to N4:N17 -----
$ws.Range("N10").Formula = "=C10-D10*20-E10*0.8-F10*0.6-H10*5+I10*10+J10/300"
$ws.Range("N11").Formula = "=C11-D11*20-E11*0.8-F11*0.6-H11*5+I11*10+J11/300"
$ws.Range("N12").Formula = "=C12-D12*20-E12*0.8-F12*0.6-H12*5+I12*10+J12/300"
$ws.Range("N13").Formula = "=C13-D13*20-E13*0.8-F13*0.6-H13*5+I13*10+J13/300"
$ws.Range("N14").Formula = "=C14-D14*20-E14*0.8-F14*0.6-H14*5+I14*10+J14/300"
$ws.Range("N15").Formula = "=C15-D15*20-E15*0.8-F15*0.6-H15*5+I15*10+J15/300"
$ws.Range("N16").Formula = "=C16-D16*20-E16*0.8-F16*0.6-H16*5+I16*10+J16/300"
$ws.Range("N17").Formula = "=C17-D17*20-E17*0.8-F17*0.6-H17*5+I17*10+J17/300"

# --- Column B width ------------------------------------------------------
$ws.Columns("B:B").ColumnWidth = 33.85546875

# --- Selection -------------------------------------------------------------
$ws.Range("G20").Select()
